$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Ishan Kishan"

# Insert a new column at A, shifting existing columns (A:L) to (B:M)
$ws.Columns.Item(1).Insert()

# Force the whole target range to be stored as text, matching the
# source data (every cell in the sheet is typed as a string in the workbook).
$ws.Range("A1:M11").NumberFormat = "@"

# Header row
$ws.Cells.Item(1, 1).Value = "matchNo"
$ws.Cells.Item(1, 2).Value = "teamName"
$ws.Cells.Item(1, 3).Value = "batterName"
$ws.Cells.Item(1, 4).Value = "states"
$ws.Cells.Item(1, 5).Value = "runs"
$ws.Cells.Item(1, 6).Value = "balls"
$ws.Cells.Item(1, 7).Value = "fours"
$ws.Cells.Item(1, 8).Value = "sixes"
$ws.Cells.Item(1, 9).Value = "sr"
$ws.Cells.Item(1, 10).Value = "opponentTeamName"
$ws.Cells.Item(1, 11).Value = "venue"
$ws.Cells.Item(1, 12).Value = "date"
$ws.Cells.Item(1, 13).Value = "result"

# Row 2
$ws.Cells.Item(2, 1).Value = "51st"
$ws.Cells.Item(2, 2).Value = "Mumbai Indians"
$ws.Cells.Item(2, 3).Value = "Ishan Kishan"
$ws.Cells.Item(2, 4).Value = ""
$ws.Cells.Item(2, 5).Value = "50"
$ws.Cells.Item(2, 6).Value = "25"
$ws.Cells.Item(2, 7).Value = "5"
$ws.Cells.Item(2, 8).Value = "3"
$ws.Cells.Item(2, 9).Value = "200.00"
$ws.Cells.Item(2, 10).Value = "Rajasthan Royals"
$ws.Cells.Item(2, 11).Value = "Sharjah"
$ws.Cells.Item(2, 12).Value = "October 05"
$ws.Cells.Item(2, 13).Value = "Mumbai won by 8 wickets (with 70 balls remaining)"

# Row 3
$ws.Cells.Item(3, 1).Value = "55th"
$ws.Cells.Item(3, 2).Value = "Mumbai Indians"
$ws.Cells.Item(3, 3).Value = "Ishan Kishan"
$ws.Cells.Item(3, 4).Value = "c †Saha b Umran Malik"
$ws.Cells.Item(3, 5).Value = "84"
$ws.Cells.Item(3, 6).Value = "32"
$ws.Cells.Item(3, 7).Value = "11"
$ws.Cells.Item(3, 8).Value = "4"
$ws.Cells.Item(3, 9).Value = "262.50"
$ws.Cells.Item(3, 10).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(3, 11).Value = "Abu Dhabi"
$ws.Cells.Item(3, 12).Value = "October 08"
$ws.Cells.Item(3, 13).Value = "Mumbai won by 42 runs"

# Row 4
$ws.Cells.Item(4, 1).Value = "13th"
$ws.Cells.Item(4, 2).Value = "Mumbai Indians"
$ws.Cells.Item(4, 3).Value = "Ishan Kishan"
$ws.Cells.Item(4, 4).Value = "b Mishra"
$ws.Cells.Item(4, 5).Value = "26"
$ws.Cells.Item(4, 6).Value = "28"
$ws.Cells.Item(4, 7).Value = "1"
$ws.Cells.Item(4, 8).Value = "1"
$ws.Cells.Item(4, 9).Value = "92.85"
$ws.Cells.Item(4, 10).Value = "Delhi Capitals"
$ws.Cells.Item(4, 11).Value = "Chennai"
$ws.Cells.Item(4, 12).Value = "April 20"
$ws.Cells.Item(4, 13).Value = "Capitals won by 6 wickets (with 5 balls remaining)"

# Row 5
$ws.Cells.Item(5, 1).Value = "5th"
$ws.Cells.Item(5, 2).Value = "Mumbai Indians"
$ws.Cells.Item(5, 3).Value = "Ishan Kishan"
$ws.Cells.Item(5, 4).Value = "c Prasidh Krishna b Cummins"
$ws.Cells.Item(5, 5).Value = "1"
$ws.Cells.Item(5, 6).Value = "3"
$ws.Cells.Item(5, 7).Value = "0"
$ws.Cells.Item(5, 8).Value = "0"
$ws.Cells.Item(5, 9).Value = "33.33"
$ws.Cells.Item(5, 10).Value = "Kolkata Knight Riders"
$ws.Cells.Item(5, 11).Value = "Chennai"
$ws.Cells.Item(5, 12).Value = "April 13"
$ws.Cells.Item(5, 13).Value = "Mumbai won by 10 runs"

# Row 6
$ws.Cells.Item(6, 1).Value = "30th"
$ws.Cells.Item(6, 2).Value = "Mumbai Indians"
$ws.Cells.Item(6, 3).Value = "Ishan Kishan"
$ws.Cells.Item(6, 4).Value = "c Raina b Bravo"
$ws.Cells.Item(6, 5).Value = "11"
$ws.Cells.Item(6, 6).Value = "10"
$ws.Cells.Item(6, 7).Value = "1"
$ws.Cells.Item(6, 8).Value = "0"
$ws.Cells.Item(6, 9).Value = "110.00"
$ws.Cells.Item(6, 10).Value = "Chennai Super Kings"
$ws.Cells.Item(6, 11).Value = "Dubai (DSC)"
$ws.Cells.Item(6, 12).Value = "September 19"
$ws.Cells.Item(6, 13).Value = "Super Kings won by 20 runs"

# Row 7
$ws.Cells.Item(7, 1).Value = "9th"
$ws.Cells.Item(7, 2).Value = "Mumbai Indians"
$ws.Cells.Item(7, 3).Value = "Ishan Kishan"
$ws.Cells.Item(7, 4).Value = "c †Bairstow b Mujeeb Ur Rahman"
$ws.Cells.Item(7, 5).Value = "12"
$ws.Cells.Item(7, 6).Value = "21"
$ws.Cells.Item(7, 7).Value = "0"
$ws.Cells.Item(7, 8).Value = "0"
$ws.Cells.Item(7, 9).Value = "57.14"
$ws.Cells.Item(7, 10).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(7, 11).Value = "Chennai"
$ws.Cells.Item(7, 12).Value = "April 17"
$ws.Cells.Item(7, 13).Value = "Mumbai won by 13 runs"

# Row 8
$ws.Cells.Item(8, 1).Value = "39th"
$ws.Cells.Item(8, 2).Value = "Mumbai Indians"
$ws.Cells.Item(8, 3).Value = "Ishan Kishan"
$ws.Cells.Item(8, 4).Value = "c Patel b Chahal"
$ws.Cells.Item(8, 5).Value = "9"
$ws.Cells.Item(8, 6).Value = "12"
$ws.Cells.Item(8, 7).Value = "1"
$ws.Cells.Item(8, 8).Value = "0"
$ws.Cells.Item(8, 9).Value = "75.00"
$ws.Cells.Item(8, 10).Value = "Royal Challengers Bangalore"
$ws.Cells.Item(8, 11).Value = "Dubai (DSC)"
$ws.Cells.Item(8, 12).Value = "September 26"
$ws.Cells.Item(8, 13).Value = "RCB won by 54 runs"

# Row 9
$ws.Cells.Item(9, 1).Value = "17th"
$ws.Cells.Item(9, 2).Value = "Mumbai Indians"
$ws.Cells.Item(9, 3).Value = "Ishan Kishan"
$ws.Cells.Item(9, 4).Value = "c †Rahul b Ravi Bishnoi"
$ws.Cells.Item(9, 5).Value = "6"
$ws.Cells.Item(9, 6).Value = "17"
$ws.Cells.Item(9, 7).Value = "0"
$ws.Cells.Item(9, 8).Value = "0"
$ws.Cells.Item(9, 9).Value = "35.29"
$ws.Cells.Item(9, 10).Value = "Punjab Kings"
$ws.Cells.Item(9, 11).Value = "Chennai"
$ws.Cells.Item(9, 12).Value = "April 23"
$ws.Cells.Item(9, 13).Value = "Punjab Kings won by 9 wickets (with 14 balls remaining)"

# Row 10
$ws.Cells.Item(10, 1).Value = "1st"
$ws.Cells.Item(10, 2).Value = "Mumbai Indians"
$ws.Cells.Item(10, 3).Value = "Ishan Kishan"
$ws.Cells.Item(10, 4).Value = "lbw b Patel"
$ws.Cells.Item(10, 5).Value = "28"
$ws.Cells.Item(10, 6).Value = "19"
$ws.Cells.Item(10, 7).Value = "2"
$ws.Cells.Item(10, 8).Value = "1"
$ws.Cells.Item(10, 9).Value = "147.36"
$ws.Cells.Item(10, 10).Value = "Royal Challengers Bangalore"
$ws.Cells.Item(10, 11).Value = "Chennai"
$ws.Cells.Item(10, 12).Value = "April 09"
$ws.Cells.Item(10, 13).Value = "RCB won by 2 wickets"

# Row 11
$ws.Cells.Item(11, 1).Value = "34th"
$ws.Cells.Item(11, 2).Value = "Mumbai Indians"
$ws.Cells.Item(11, 3).Value = "Ishan Kishan"
$ws.Cells.Item(11, 4).Value = "c Russell b Ferguson"
$ws.Cells.Item(11, 5).Value = "14"
$ws.Cells.Item(11, 6).Value = "13"
$ws.Cells.Item(11, 7).Value = "0"
$ws.Cells.Item(11, 8).Value = "1"
$ws.Cells.Item(11, 9).Value = "107.69"
$ws.Cells.Item(11, 10).Value = "Kolkata Knight Riders"
$ws.Cells.Item(11, 11).Value = "Abu Dhabi"
$ws.Cells.Item(11, 12).Value = "September 23"
$ws.Cells.Item(11, 13).Value = "KKR won by 7 wickets (with 29 balls remaining)"

